$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset counts (materias, profesores, estudiantes) to 0 for rows 2-4
$ws.Range("B2:D4").Value = 0
